{"js": "// Update the answer values in the \"three-digit x one-digit\" practice\n// table. Every populated row (0-based row indices 0, 4, 9, 14, 19) gets\n// its five cell values replaced with a new set of multiplication facts.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst rowUpdates = {\n  0: [\"165\u00d77=1155\", \"891\u00d76=5346\", \"493\u00d76=2958\", \"566\u00d79=5094\", \"424\u00d76=2544\"],\n  4: [\"341\u00d74=1364\", \"923\u00d77=6461\", \"997\u00d74=3988\", \"101\u00d79=909\", \"226\u00d74=904\"],\n  9: [\"446\u00d76=2676\", \"790\u00d74=3160\", \"593\u00d77=4151\", \"224\u00d72=448\", \"633\u00d79=5697\"],\n  14: [\"875\u00d75=4375\", \"435\u00d79=3915\", \"740\u00d79=6660\", \"789\u00d76=4734\", \"254\u00d77=1778\"],\n  19: [\"849\u00d77=5943\", \"170\u00d73=510\", \"162\u00d73=486\", \"448\u00d79=4032\", \"344\u00d75=1720\"],\n};\n\nfor (const rowIndex of Object.keys(rowUpdates)) {\n  const values = rowUpdates[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(Number(rowIndex), col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# Map of 1-based row index -> new cell texts (left to right) for the rows\n# that contain the \"###x#=####\" practice problems.\n$rowsData = @{\n    1  = @('165\u00d77=1155', '891\u00d76=5346', '493\u00d76=2958', '566\u00d79=5094', '424\u00d76=2544')\n    5  = @('341\u00d74=1364', '923\u00d77=6461', '997\u00d74=3988', '101\u00d79=909',  '226\u00d74=904')\n    10 = @('446\u00d76=2676', '790\u00d74=3160', '593\u00d77=4151', '224\u00d72=448',  '633\u00d79=5697')\n    15 = @('875\u00d75=4375', '435\u00d79=3915', '740\u00d79=6660', '789\u00d76=4734', '254\u00d77=1778')\n    20 = @('849\u00d77=5943', '170\u00d73=510',  '162\u00d73=486',  '448\u00d79=4032', '344\u00d75=1720')\n}\n\nforeach ($rowIndex in $rowsData.Keys) {\n    $values = $rowsData[$rowIndex]\n    $row = $t.Rows($rowIndex)\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $cell = $row.Cells($col)\n        $cell.Range.Text = $values[$col - 1]\n    }\n}\n"}
